# Fruta / hortaliza, semanal
# Insert a new weekly data block (3 rows: Especial/Primera/Segunda) for
# "Comercializadora del Agro de Limari" right before the existing block
# that starts at row 207, pushing the remaining rows (207:274) down to
# (210:277) and extending the used range to A1:T277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 207 - everything currently at 207:274 shifts to 210:277
$ws.Rows("207:209").Insert()

# Shared/constant values for this data block (match the surrounding rows)
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100101
$producto    = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "`$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7
$fecha       = 44524

# New rows data: row, calidad, volumen, precioMin, precioMax, precioProm, precioKg
$rows = @(
    @{ R = 207; Calidad = "Especial"; Volumen = 400; PMin = 12500; PMax = 13000; PProm = 12750; PKg = 1821 },
    @{ R = 208; Calidad = "Primera";  Volumen = 300; PMin = 10500; PMax = 11000; PProm = 10750; PKg = 1536 },
    @{ R = 209; Calidad = "Segunda";  Volumen = 200; PMin = 8500;  PMax = 9000;  PProm = 8750;  PKg = 1250 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $row.Calidad
    $ws.Cells.Item($r, 13).Value = $row.Volumen
    $ws.Cells.Item($r, 14).Value = $row.PMin
    $ws.Cells.Item($r, 15).Value = $row.PMax
    $ws.Cells.Item($r, 16).Value = $row.PProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $row.PKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
